$wb = $excel.ActiveWorkbook

# Rename "Sheet5" -> "HodoCalibRuns" (it holds the Hodoscope calibration run table)
$hodoSheet = $wb.Worksheets.Item("Sheet5")
$hodoSheet.Name = "HodoCalibRuns"

# --- Sheet1: scroll back to top, move selection to L24 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("L24").Select() | Out-Null

# --- MissinginMaurik: scroll back to top (col A), move selection to L19 ---
$ws2 = $wb.Worksheets.Item("MissinginMaurik")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("L19").Select() | Out-Null

# --- MissingCameron: scroll back to top, keep selection on L27 ---
$ws3 = $wb.Worksheets.Item("MissingCameron")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("L27").Select() | Out-Null

# --- HodoCalibRuns: scroll to row 3, move selection to H14, widen new columns E/F ---
$ws5 = $wb.Worksheets.Item("HodoCalibRuns")
$ws5.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws5.Range("H14").Select() | Out-Null

# New/adjusted column widths (characters) for columns E and F
$ws5.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws5.Columns.Item(6).ColumnWidth = 11.333333333333334

# Leave HodoCalibRuns as the active/selected tab, matching tabSelected="true"
$ws5.Activate()
